$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Change 1: insert a new "Meta description" paragraph right after the
# first (Heading1) paragraph.
# -----------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
[void]$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:r/>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">: Abby and The Witch is an engaging slot game with visually engaging graphics, free spins mode, and low volatility. Read our review and play for free.</w:t></w:r>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

[void]$metaPara.Range.InsertXML($metaXml)

# -----------------------------------------------------------------------
# Change 2: drop the trailing bold "Play Abby and The Witch for Free -
# Slot Game Review" paragraph near the end of the document, and update
# the italic paragraph that follows it with the new image prompt text.
# -----------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldPara = $d.Paragraphs($count - 1)
$italicPara = $d.Paragraphs($count)

if ($boldPara.Range.Text.Trim() -ne "Play Abby and The Witch for Free - Slot Game Review") {
    throw "Unexpected paragraph content before deletion: $($boldPara.Range.Text)"
}

$delRange = $d.Range($boldPara.Range.Start, $italicPara.Range.Start)
[void]$delRange.Delete()

$count = $d.Paragraphs.Count
$italicPara = $d.Paragraphs($count)
$pr = $italicPara.Range
$textRange = $d.Range($pr.Start, $pr.End - 1)
$textRange.Text = "Prompt: Create a feature image for Abby & The Witch Design a cartoon-style feature image that includes a happy-looking Maya warrior wearing glasses. The image should also incorporate elements from the game `"Abby & The Witch,`" such as Abby herself, the colorless world, and Baba Yaga's house and cemetery. Use bright colors to contrast the black and white world of the game and make the Maya warrior stand out. Feel free to add other magical elements to the image, like spells, potions, or magical creatures, to give it a more whimsical feel. The image should be eye-catching and convey the spirit of adventure and magic that the game offers to players."
